$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new column L (2021 data) is added, mirroring column K (2020 data).
# First, copy the formatting of K3:K11 into L3:L11 so the new column matches
# the look (borders, fonts, number formats) of the last existing year column.
for ($r = 3; $r -le 11; $r++) {
    $ws.Cells.Item($r, 11).Copy()
    $ws.Cells.Item($r, 12).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

# Now fill in the values for column L (row 3 stays blank, like K3)
$ws.Range("L4").Value = 2021
$ws.Range("L5").Value = 0.86
$ws.Range("L6").Value = 1.07
$ws.Range("L7").Value = 25.27
$ws.Range("L8").Value = 14
$ws.Range("L9").Value = 0.12
$ws.Range("L10").Value = 21.74
$ws.Range("L11").Value = 9.4600000000000009

# Update the selected cell shown when the workbook was last saved
$ws.Range("N2").Select()
